# Calibration run with ATB 2023 data that removes 111 rules
#
# The calibrated "share of costs that must be covered" values on the
# SoCtMbCtbDP sheet move from 0.95 to 1 for every data row (B2:B25), and
# that sheet becomes the active/selected tab (it was "About" before).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SoCtMbCtbDP")

# Recalibrate every data row: the share of costs that must be covered to
# be deemed profitable is now 100% (1) instead of 95% (0.95).
$ws.Range("B2:B25").Value = 1

# Make SoCtMbCtbDP the active sheet/tab (was "About").
$ws.Activate()
$ws.Range("B2:B25").Select()
